$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New attendance rows: day 20 (row 25 gets Day-Count/Date added) plus
# six more rows (26-31) through 26-02-2026.

# Cells that need to be added, keyed by A1 ref -> value
$newCells = @{
    "B25" = "21"
    "C25" = "20-02-20206"

    "A26" = "25"
    "B26" = "22"
    "C26" = "21-02-2026"

    "A27" = "26"
    "B27" = "null"
    "C27" = "22-02-2026"
    "D27" = "sunday"

    "A28" = "27"
    "B28" = "23"
    "C28" = "23-02-2026"

    "A29" = "28"
    "B29" = "24"
    "C29" = "24-02-2026"

    "A30" = "29"
    "B30" = "25"
    "C30" = "25-02-2026"

    "A31" = "30"
    "B31" = "26"
    "C31" = "26-02-2026"
}

# Force the destination cells to be stored as text (matching the rest of
# the sheet, where numeric-looking values are kept as text / flagged via
# numberStoredAsText) while keeping the same "General" number format /
# style as the existing rows once we flip the format back. NumberFormat
# is applied per-cell (a comma-unioned multi-area Range only honors the
# first area for this engine), so loop individually.
$refs = [string[]]$newCells.Keys

foreach ($ref in $refs) {
    $ws.Range($ref).NumberFormat = "@"
    $ws.Range($ref).Value = $newCells[$ref]
    $ws.Range($ref).NumberFormat = "General"
}

# Keep the same visual style as the rest of the table (style index 1).
$ws.Range("A25").Copy() | Out-Null
foreach ($ref in $refs) {
    $ws.Range($ref).PasteSpecial(-4122) | Out-Null
}
